# Insert a new data row at row 65 (pushing the existing row 65..168 down to
# 66..169), then populate the new row with its own data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value = 4
$ws.Cells.Item(65, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(65, 3).Value = "Los Lagos"
$ws.Cells.Item(65, 4).Value = 44477
$ws.Cells.Item(65, 5).Value = 10
$ws.Cells.Item(65, 6).Value = 100112037
$ws.Cells.Item(65, 7).Value = "Cebollín"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 160
$ws.Cells.Item(65, 11).Value = 6000
$ws.Cells.Item(65, 12).Value = 6000
$ws.Cells.Item(65, 13).Value = 6000
$ws.Cells.Item(65, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(65, 15).Value = "Región Metropolitana"
$ws.Cells.Item(65, 16).Value = 167
$ws.Cells.Item(65, 17).Value = 36
$ws.Cells.Item(65, 18).Value = "Hortaliza"
